$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Update the date line: "27 Feb 2021 17:59:29" -> "2 Mar 2021 19:38:53"
# ------------------------------------------------------------------
$d.Content.Find.Execute("27", $true, $true, $false, $false, $false, $true, 1, $false, "2", 2)
$d.Content.Find.Execute("Feb", $true, $true, $false, $false, $false, $true, 1, $false, "Mar", 2)
$d.Content.Find.Execute("17:59:29", $true, $true, $false, $false, $false, $true, 1, $false, "19:38:53", 2)

# ------------------------------------------------------------------
# 2. Update the "use" command to point to the full URL for the data file.
#    Note: we avoid Find.Execute's built-in Replace for strings containing
#    double quotes because it triggers smart-quote autocorrection; instead
#    locate the text with Find and then assign Range.Text directly.
# ------------------------------------------------------------------
$r = $d.Content
if ($r.Find.Execute(". use Spruce.dta, clear")) {
  $r.Text = ". use `"https://github.com/agrogan1/multilevel/raw/master/reshaping-data/Spruce.dta`", clear"
}

# ------------------------------------------------------------------
# 3. Update the "Contains data from Spruce.dta" line to reference the new URL
# ------------------------------------------------------------------
$d.Content.Find.Execute("Contains data from Spruce.dta", $true, $true, $false, $false, $false, $true, 1, $false, "Contains data from https://github.com/agrogan1/multilevel/raw/master/reshaping-data/Spruce.dta", 2)

# ------------------------------------------------------------------
# 4. Insert a new "Describe The Data" Heading1 paragraph (with bookmark)
#    right before the ". describe" source code block, i.e. right after
#    the paragraph containing the "use ..." command.
# ------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs.Item($i)
  if ($p.Range.Text -like "*use *Spruce.dta*") {
    $target = $p
  }
}

$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:bookmarkStart w:id="23" w:name="describe-the-data"/><w:bookmarkEnd w:id="23"/><w:r><w:t xml:space="preserve">Describe The Data 🌲</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.Range.InsertParagraphAfter()
$newPara = $target.Next()
$newPara.Range.InsertXML($newParaXml)
